# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The sheet "Hoja1" holds a worker/period table in B15:J26 (B15:J15 are the
# headers). Rows 16-26 are re-populated here: the records were reordered
# (grouped back together per worker) and the GUSTAVO HUMBERTO VANEGAS GARCIA
# "Salario Basico" (column G) values were corrected from 1200000 to 1000000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row=16; B="CC"; C="49724403";   D="CARMEN ISABEL MANJARRES CARCAMO"; E="1812"; F=31249; G=781242 },
    @{ Row=17; B="CC"; C="1047439413"; D="GUSTAVO HUMBERTO VANEGAS GARCIA"; E="1812"; F=40000; G=1000000 },
    @{ Row=18; B="CC"; C="1050952506"; D="YURIS ROCIO PUELLO OYOLA";        E="1903"; F=40000; G=1000000 },
    @{ Row=19; B="CC"; C="73183995";   D="STALIN RODRIGO CHAPUEL TELLO";    E="1903"; F=33125; G=828116 },
    @{ Row=20; B="CC"; C="1047439413"; D="GUSTAVO HUMBERTO VANEGAS GARCIA"; E="1903"; F=48000; G=1000000 },
    @{ Row=21; B="CC"; C="1050952506"; D="YURIS ROCIO PUELLO OYOLA";        E="1904"; F=40000; G=1000000 },
    @{ Row=22; B="CC"; C="73183995";   D="STALIN RODRIGO CHAPUEL TELLO";    E="1904"; F=33125; G=828116 },
    @{ Row=23; B="CC"; C="1047439413"; D="GUSTAVO HUMBERTO VANEGAS GARCIA"; E="1904"; F=40000; G=1000000 },
    @{ Row=24; B="CC"; C="73183995";   D="STALIN RODRIGO CHAPUEL TELLO";    E="1905"; F=33125; G=828116 },
    @{ Row=25; B="CC"; C="73183995";   D="STALIN RODRIGO CHAPUEL TELLO";    E="1906"; F=33125; G=828116 },
    @{ Row=26; B="CC"; C="73183995";   D="STALIN RODRIGO CHAPUEL TELLO";    E="1907"; F=16562; G=828116 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
